$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.986.54'
$ws.Range('D3').Value = '3.120.61'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('E9').Value = '  -3.15%  '
$ws.Range('E10').Value = '  -1.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.481'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = '3.637.49'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '66.937.81'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '3.119.03'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '475.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.708'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '83.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('E25').Value = '  -3.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.83%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.91'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.85%  '
$ws.Range('E29').Value = '  -1.42%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.57'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.114'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').Value = '0.0₃0949'
$ws.Range('E33').Value = '  -6.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.83'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.976'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '46.96'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.33%  '
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.05'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.311'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('E41').Value = '  +1.05%  '
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('D43').Value = '2.822.95'
$ws.Range('E43').Value = '  +0.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '382.67'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('E46').Value = '  -9.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '135.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('E50').Value = '  -1.19%  '
$ws.Range('E51').Value = '  -0.87%  '
